# Refresh the stock-screener table (Sheet1) with a new set of tickers
# in columns B (Buying Opportunity), C (support Zone), D (long buildup)
# and E (Short buildup); column F (FII ENTERING) is cleared, and the
# trailing rows that are no longer needed are removed so the used range
# shrinks from A1:F34 down to A1:F29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column index map: A=1 B=2 C=3 D=4 E=5 F=6

$ws.Cells.Item(2, 2).Value = "NSE:BBOX"
$ws.Cells.Item(2, 3).Value = "NSE:ACC"
$ws.Cells.Item(2, 4).Value = "NSE:KOTAKBANK"
$ws.Cells.Item(2, 5).Value = "NSE:ASTRAL"
$ws.Cells.Item(2, 6).Value = ""

$ws.Cells.Item(3, 2).Value = "NSE:CUBEXTUB"
$ws.Cells.Item(3, 3).Value = "NSE:ADFFOODS"
$ws.Cells.Item(3, 4).Value = ""
$ws.Cells.Item(3, 5).Value = "NSE:BAJAJFINSV"
$ws.Cells.Item(3, 6).Value = ""

$ws.Cells.Item(4, 2).Value = "NSE:DIVGIITTS"
$ws.Cells.Item(4, 3).Value = "NSE:BALKRISIND"
$ws.Cells.Item(4, 4).Value = ""
$ws.Cells.Item(4, 5).Value = "NSE:COLPAL"
$ws.Cells.Item(4, 6).Value = ""

$ws.Cells.Item(5, 2).Value = "NSE:ENGINERSIN"
$ws.Cells.Item(5, 3).Value = "NSE:BOROLTD"
$ws.Cells.Item(5, 4).Value = ""
$ws.Cells.Item(5, 5).Value = "NSE:DABUR"
$ws.Cells.Item(5, 6).Value = ""

$ws.Cells.Item(6, 2).Value = "NSE:GICHSGFIN"
$ws.Cells.Item(6, 3).Value = "NSE:CANTABIL"
$ws.Cells.Item(6, 4).Value = ""
$ws.Cells.Item(6, 5).Value = "NSE:DIXON"
$ws.Cells.Item(6, 6).Value = ""

$ws.Cells.Item(7, 2).Value = "NSE:GMBREW"
$ws.Cells.Item(7, 3).Value = "NSE:DABUR"
$ws.Cells.Item(7, 4).Value = ""
$ws.Cells.Item(7, 5).Value = "NSE:HDFCAMC"
$ws.Cells.Item(7, 6).Value = ""

$ws.Cells.Item(8, 2).Value = "NSE:GMDCLTD"
$ws.Cells.Item(8, 3).Value = "NSE:DCXINDIA"
$ws.Cells.Item(8, 4).Value = ""
$ws.Cells.Item(8, 5).Value = "NSE:IOC"
$ws.Cells.Item(8, 6).Value = ""

$ws.Cells.Item(9, 2).Value = "NSE:GRAVITA"
$ws.Cells.Item(9, 3).Value = "NSE:ELGIEQUIP"
$ws.Cells.Item(9, 4).Value = ""
$ws.Cells.Item(9, 5).Value = "NSE:KAYNES"
$ws.Cells.Item(9, 6).Value = ""

$ws.Cells.Item(10, 2).Value = "NSE:HDFCGOLD"
$ws.Cells.Item(10, 3).Value = "NSE:FINOPB"
$ws.Cells.Item(10, 4).Value = ""
$ws.Cells.Item(10, 5).Value = "NSE:LICI"
$ws.Cells.Item(10, 6).Value = ""

$ws.Cells.Item(11, 2).Value = "NSE:HDFCSILVER"
$ws.Cells.Item(11, 3).Value = "NSE:FSL"
$ws.Cells.Item(11, 4).Value = ""
$ws.Cells.Item(11, 5).Value = "NSE:LT"
$ws.Cells.Item(11, 6).Value = ""

$ws.Cells.Item(12, 2).Value = "NSE:HINDCOPPER"
$ws.Cells.Item(12, 3).Value = "NSE:GLAXO"
$ws.Cells.Item(12, 4).Value = ""
$ws.Cells.Item(12, 5).Value = "NSE:MFSL"
$ws.Cells.Item(12, 6).Value = ""

$ws.Cells.Item(13, 2).Value = "NSE:MAITHANALL"
$ws.Cells.Item(13, 3).Value = "NSE:GRWRHITECH"
$ws.Cells.Item(13, 4).Value = ""
$ws.Cells.Item(13, 5).Value = ""
$ws.Cells.Item(13, 6).Value = ""

$ws.Cells.Item(14, 2).Value = "NSE:MOIL"
$ws.Cells.Item(14, 3).Value = "NSE:HINDPETRO"
$ws.Cells.Item(14, 4).Value = ""
$ws.Cells.Item(14, 5).Value = ""
$ws.Cells.Item(14, 6).Value = ""

$ws.Cells.Item(15, 2).Value = "NSE:PARACABLES"
$ws.Cells.Item(15, 3).Value = "NSE:ICICIBANK"
$ws.Cells.Item(15, 4).Value = ""
$ws.Cells.Item(15, 5).Value = ""
$ws.Cells.Item(15, 6).Value = ""

$ws.Cells.Item(16, 2).Value = "NSE:RAIN"
$ws.Cells.Item(16, 3).Value = "NSE:INFY"
$ws.Cells.Item(16, 4).Value = ""
$ws.Cells.Item(16, 5).Value = ""
$ws.Cells.Item(16, 6).Value = ""

$ws.Cells.Item(17, 2).Value = ""
$ws.Cells.Item(17, 3).Value = "NSE:IZMO"
$ws.Cells.Item(17, 4).Value = ""
$ws.Cells.Item(17, 5).Value = ""
$ws.Cells.Item(17, 6).Value = ""

$ws.Cells.Item(18, 2).Value = ""
$ws.Cells.Item(18, 3).Value = "NSE:KEI"
$ws.Cells.Item(18, 4).Value = ""
$ws.Cells.Item(18, 5).Value = ""
$ws.Cells.Item(18, 6).Value = ""

$ws.Cells.Item(19, 2).Value = ""
$ws.Cells.Item(19, 3).Value = "NSE:LATENTVIEW"
$ws.Cells.Item(19, 4).Value = ""
$ws.Cells.Item(19, 5).Value = ""
$ws.Cells.Item(19, 6).Value = ""

$ws.Cells.Item(20, 2).Value = ""
$ws.Cells.Item(20, 3).Value = "NSE:LT"
$ws.Cells.Item(20, 4).Value = ""
$ws.Cells.Item(20, 5).Value = ""
$ws.Cells.Item(20, 6).Value = ""

$ws.Cells.Item(21, 2).Value = ""
$ws.Cells.Item(21, 3).Value = "NSE:MANALIPETC"
$ws.Cells.Item(21, 4).Value = ""
$ws.Cells.Item(21, 5).Value = ""
$ws.Cells.Item(21, 6).Value = ""

$ws.Cells.Item(22, 2).Value = ""
$ws.Cells.Item(22, 3).Value = "NSE:MAXIND"
$ws.Cells.Item(22, 4).Value = ""
$ws.Cells.Item(22, 5).Value = ""
$ws.Cells.Item(22, 6).Value = ""

$ws.Cells.Item(23, 2).Value = ""
$ws.Cells.Item(23, 3).Value = "NSE:MSUMI"
$ws.Cells.Item(23, 4).Value = ""
$ws.Cells.Item(23, 5).Value = ""
$ws.Cells.Item(23, 6).Value = ""

$ws.Cells.Item(24, 2).Value = ""
$ws.Cells.Item(24, 3).Value = "NSE:NEULANDLAB"
$ws.Cells.Item(24, 4).Value = ""
$ws.Cells.Item(24, 5).Value = ""
$ws.Cells.Item(24, 6).Value = ""

$ws.Cells.Item(25, 2).Value = ""
$ws.Cells.Item(25, 3).Value = "NSE:PAGEIND"
$ws.Cells.Item(25, 4).Value = ""
$ws.Cells.Item(25, 5).Value = ""
$ws.Cells.Item(25, 6).Value = ""

$ws.Cells.Item(26, 2).Value = ""
$ws.Cells.Item(26, 3).Value = "NSE:PIONEEREMB"
$ws.Cells.Item(26, 4).Value = ""
$ws.Cells.Item(26, 5).Value = ""
$ws.Cells.Item(26, 6).Value = ""

$ws.Cells.Item(27, 2).Value = ""
$ws.Cells.Item(27, 3).Value = "NSE:PRAJIND"
$ws.Cells.Item(27, 4).Value = ""
$ws.Cells.Item(27, 5).Value = ""
$ws.Cells.Item(27, 6).Value = ""

$ws.Cells.Item(28, 2).Value = ""
$ws.Cells.Item(28, 3).Value = "NSE:PRUDENT"
$ws.Cells.Item(28, 4).Value = ""
$ws.Cells.Item(28, 5).Value = ""
$ws.Cells.Item(28, 6).Value = ""

$ws.Cells.Item(29, 2).Value = ""
$ws.Cells.Item(29, 3).Value = "NSE:RRKABEL"
$ws.Cells.Item(29, 4).Value = ""
$ws.Cells.Item(29, 5).Value = ""
$ws.Cells.Item(29, 6).Value = ""

# Remove now-unused rows 30-34 (previously rows 28-32 of stock data)
$ws.Rows("30:34").Delete()
